$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value corrections -------------------------------------------------
$ws.Range("K3").Value = 2144
$ws.Range("Q10").Value = 29
$ws.Range("G13").Value = 2183

# --- Insert a new manufacturer row ("Viking") before the "Vulcanair" row -----
# Before the edit: row 69 = Vulcanair, 70 = Wassmer, 71 = Zenair.
# After the edit:  row 69 = Viking (new), 70 = Vulcanair, 71 = Wassmer, 72 = Zenair.
$ws.Rows.Item(69).Insert()

# Copy formatting (border/font/alignment) from the row below (now Vulcanair)
# onto the freshly inserted, blank row so the new "Viking" row matches the
# look of every other data row.
$ws.Range("A70").Copy()
$ws.Range("A69").PasteSpecial(-4122)

# Populate the new row's data: all aircraft-type counts are 0 except "Other".
$ws.Range("A69").Value = "Viking"
$ws.Range("B69:P69").Value = 0
$ws.Range("Q69").Value = 2
